$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '90.735.62'
$ws.Cells.Item(2, 5).Value = '  +0.28%  '

$ws.Cells.Item(3, 4).Value = '3.115.69'
$ws.Cells.Item(3, 5).Value = '  +0.21%  '

$ws.Cells.Item(4, 5).Value = '  -0.24%  '

$ws.Cells.Item(5, 4).Value = '''241.74'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +1.67%  '

$ws.Cells.Item(6, 4).Value = '''624.41'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.44%  '

$ws.Cells.Item(7, 5).Value = '  +15.10%  '

$ws.Cells.Item(8, 5).Value = '  +5.15%  '

$ws.Cells.Item(9, 5).Value = '  +0.10%  '

$ws.Cells.Item(10, 4).Value = '3.111.56'
$ws.Cells.Item(10, 5).Value = '  +0.11%  '

$ws.Cells.Item(11, 4).Value = '''0.767'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +7.86%  '

$ws.Cells.Item(12, 4).Value = '''0.204'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +4.06%  '

$ws.Cells.Item(13, 4).Value = '''0.0000251'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +4.13%  '

$ws.Cells.Item(14, 4).Value = '''35.33'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -3.10%  '

$ws.Cells.Item(15, 4).Value = '''5.50'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -1.34%  '

$ws.Cells.Item(16, 4).Value = '90.549.65'
$ws.Cells.Item(16, 5).Value = '  +0.41%  '

$ws.Cells.Item(17, 4).Value = '3.685.64'
$ws.Cells.Item(17, 5).Value = '  +0.23%  '

$ws.Cells.Item(18, 4).Value = '3.093.95'
$ws.Cells.Item(18, 5).Value = '  +0.11%  '

$ws.Cells.Item(19, 5).Value = '  +2.81%  '

$ws.Cells.Item(20, 4).Value = '''14.45'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +1.36%  '

$ws.Cells.Item(21, 2).Value = 'Polkadot'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(21, 4).Value = '''5.89'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +5.10%  '

$ws.Cells.Item(22, 2).Value = 'PEPE'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(22, 4).Value = '''0.0000210'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.71%  '

$ws.Cells.Item(23, 4).Value = '''450.48'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.95%  '

$ws.Cells.Item(24, 5).Value = '  +1.82%  '

$ws.Cells.Item(25, 4).Value = '''5.92'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -2.37%  '

$ws.Cells.Item(26, 4).Value = '''93.41'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +4.64%  '

$ws.Cells.Item(27, 4).Value = '''11.87'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -3.84%  '

$ws.Cells.Item(28, 4).Value = '3.269.95'
$ws.Cells.Item(28, 5).Value = '  +0.03%  '

$ws.Cells.Item(29, 5).Value = '  +0.07%  '

$ws.Cells.Item(30, 4).Value = '''0.180'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +13.40%  '

$ws.Cells.Item(31, 4).Value = '''0.233'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +17.97%  '

$ws.Cells.Item(32, 4).Value = '''0.120'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +40.36%  '

$ws.Cells.Item(33, 4).Value = '''9.11'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -2.27%  '

$ws.Cells.Item(34, 4).Value = '''0.999'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +35.80%  '

$ws.Cells.Item(35, 4).Value = '''0.163'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +8.15%  '

$ws.Cells.Item(36, 4).Value = '''26.85'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -1.74%  '

$ws.Cells.Item(37, 4).Value = '''7.71'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +9.75%  '

$ws.Cells.Item(38, 4).Value = '''4.22'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +27.95%  '

$ws.Cells.Item(39, 2).Value = 'Bittensor'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(39, 4).Value = '''495.89'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -2.29%  '

$ws.Cells.Item(40, 2).Value = 'PancakeSwap'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(40, 4).Value = '''1.92'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +0.14%  '

$ws.Cells.Item(41, 4).Value = '''3.59'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -4.81%  '

$ws.Cells.Item(42, 4).Value = '''1.30'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -1.22%  '

$ws.Cells.Item(43, 4).Value = '''0.420'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -0.21%  '

$ws.Cells.Item(44, 4).Value = '''22.11'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -0.26%  '

$ws.Cells.Item(46, 4).Value = '''157.57'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +5.81%  '

$ws.Cells.Item(47, 4).Value = '''1.92'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -1.23%  '

$ws.Cells.Item(48, 4).Value = '''0.697'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +0.35%  '

$ws.Cells.Item(49, 4).Value = '''4.61'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +2.73%  '

$ws.Cells.Item(50, 2).Value = 'ImmutableX'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(50, 4).Value = '''1.35'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +1.30%  '

$ws.Cells.Item(51, 2).Value = 'OKB'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(51, 4).Value = '''45.03'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -0.76%  '
